# Fruta / hortaliza, semanal
# Insert a new weekly record at row 36, pushing the existing rows 36-51 down
# to 37-52 (values unchanged, just shifted), and fill the new row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 36; this shifts rows 36..51 down
# to 37..52 and copies row 36's formatting (including the date style on D)
# into the freshly inserted row.
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new observation.
$ws.Cells.Item(36, 1).Value  = 7
$ws.Cells.Item(36, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value  = "Ñuble"
$ws.Cells.Item(36, 4).Value  = 44572
$ws.Cells.Item(36, 5).Value  = 16
$ws.Cells.Item(36, 6).Value  = 100112022
$ws.Cells.Item(36, 7).Value  = "Arveja Verde"
$ws.Cells.Item(36, 8).Value  = "Sin especificar"
$ws.Cells.Item(36, 9).Value  = "Primera"
$ws.Cells.Item(36, 10).Value = 60
$ws.Cells.Item(36, 11).Value = 27000
$ws.Cells.Item(36, 12).Value = 28000
$ws.Cells.Item(36, 13).Value = 27500
$ws.Cells.Item(36, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(36, 16).Value = 1100
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
